# Update cryptos list values (Price / Volume(1h)) per commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.562.89"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.983.42"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0846"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("D13").Value = "3.447.25"
$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("E15").Value = "  +2.83%  "

$ws.Range("D16").Value = "2.978.00"
$ws.Range("E16").Value = "  +2.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.975"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.70%  "

$ws.Range("D18").Value = "51.503.82"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +2.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +22.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.117"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.02%  "

$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.95%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.96%  "

$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("E42").Value = "  +2.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.03%  "

$ws.Range("E45").Value = "  +18.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "

$ws.Range("D48").Value = "2.034.15"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("E50").Value = "  +8.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.00%  "
